# Update the "想去人数" (number of people interested) figures in column F
# for both the "展览" and "全部类型" worksheets, which carry the same data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    2  = 8399
    3  = 7965
    4  = 132
    9  = 132
    12 = 718
    13 = 139
    14 = 2009
    19 = 133
    20 = 30
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
